# Restore C10 on the active sheet to the value 1 (was 18).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 1
